$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.27"
$ws.Range("E2").Value = "'0.82%"
$ws.Range("D3").Value = "'29.70"
$ws.Range("E3").Value = "'9.55%"
$ws.Range("D4").Value = "'5.186"
$ws.Range("E4").Value = "'3.23%"
$ws.Range("D5").Value = "'0.05702"
$ws.Range("E5").Value = "'0.67%"
$ws.Range("D6").Value = "'6.592"
$ws.Range("E6").Value = "'1.92%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8578"
$ws.Range("E7").Value = "'4.36%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.8776"
$ws.Range("E8").Value = "'4.28%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1368"
$ws.Range("E9").Value = "'3.28%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07090"
$ws.Range("E10").Value = "'2.51%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.02863"
$ws.Range("E11").Value = "'-0.69%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09382"
$ws.Range("E12").Value = "'-0.10%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001530"
$ws.Range("E13").Value = "'0.65%"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "'0.04149"
$ws.Range("E14").Value = "'0.35%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005978"
$ws.Range("E15").Value = "'-0.20%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006084"
$ws.Range("E16").Value = "'-2.08%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007491"
$ws.Range("E17").Value = "'5,107.52%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.483"
$ws.Range("E18").Value = "'-0.80%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.054"
$ws.Range("E19").Value = "'1.81%"
$ws.Range("D20").Value = "'2.276"
$ws.Range("E20").Value = "'2.21%"
$ws.Range("D21").Value = "'0.3144"
$ws.Range("D22").Value = "'0.03270"
$ws.Range("E22").Value = "'3.76%"
$ws.Range("D23").Value = "'0.1301"
$ws.Range("E23").Value = "'0.76%"
$ws.Range("D24").Value = "'2.906"
$ws.Range("E24").Value = "'-18.71%"
$ws.Range("D26").Value = "'0.005089"
$ws.Range("E26").Value = "'31.51%"
$ws.Range("D27").Value = "'0.001219"
$ws.Range("E27").Value = "'-0.11%"
$ws.Range("D28").Value = "'0.0001209"
$ws.Range("E28").Value = "'23.42%"
$ws.Range("D40").Value = "'0.03744"
$ws.Range("E40").Value = "'2.06%"
$ws.Range("D41").Value = "'0.005678"
$ws.Range("E41").Value = "'65.11%"
$ws.Range("E42").Value = "'1.62%"
$ws.Range("D43").Value = "'0.002099"
$ws.Range("E43").Value = "'-7.82%"
$ws.Range("D44").Value = "'0.009410"
$ws.Range("E44").Value = "'-17.05%"
$ws.Range("D45").Value = "'0.00005107"
$ws.Range("E45").Value = "'-3.89%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("D47").Value = "'0.07097"
$ws.Range("E47").Value = "'-30.07%"
$ws.Range("D48").Value = "'0.002671"
$ws.Range("E48").Value = "'3.91%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.03%"
